# Natmi following Dr Hou advice
# Recompute LR-pair stats (Fstl1-Dip2a) after expression input update:
# Ligand/Receptor expressing-cell counts rise from 1 to 3, and all
# dependent average/total/specificity/edge-weight columns are refreshed
# to match the updated NATMI output for every data row (rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 21.395034
$ws.Cells.Item(2, 8).Value = 64.185102
$ws.Cells.Item(2, 9).Value = 0.03113537258663955
$ws.Cells.Item(2, 10).Value = 0.03113537258663955
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 7.321725333333333
$ws.Cells.Item(2, 14).Value = 21.965176
$ws.Cells.Item(2, 15).Value = 0.2286194065168947
$ws.Cells.Item(2, 16).Value = 0.2286194065168946
$ws.Cells.Item(2, 17).Value = 156.648562445328
$ws.Cells.Item(2, 18).Value = 1409.837062007952
$ws.Cells.Item(2, 19).Value = 0.007118150402439924
$ws.Cells.Item(2, 20).Value = 0.007118150402439924
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 21.395034
$ws.Cells.Item(3, 8).Value = 64.185102
$ws.Cells.Item(3, 9).Value = 0.03113537258663955
$ws.Cells.Item(3, 10).Value = 0.03113537258663955
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 14.81128766666667
$ws.Cells.Item(3, 14).Value = 44.433863
$ws.Cells.Item(3, 15).Value = 0.4624794897301531
$ws.Cells.Item(3, 16).Value = 0.462479489730153
$ws.Cells.Item(3, 17).Value = 316.888003212114
$ws.Cells.Item(3, 18).Value = 2851.992028909026
$ws.Cells.Item(3, 19).Value = 0.01439947122642725
$ws.Cells.Item(3, 20).Value = 0.01439947122642725
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 21.395034
$ws.Cells.Item(4, 8).Value = 64.185102
$ws.Cells.Item(4, 9).Value = 0.03113537258663955
$ws.Cells.Item(4, 10).Value = 0.03113537258663955
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 2.522212
$ws.Cells.Item(4, 14).Value = 7.566636
$ws.Cells.Item(4, 15).Value = 0.078755564337357
$ws.Cells.Item(4, 16).Value = 0.07875556433735698
$ws.Cells.Item(4, 17).Value = 53.962811495208
$ws.Cells.Item(4, 18).Value = 485.665303456872
$ws.Cells.Item(4, 19).Value = 0.002452083838914672
$ws.Cells.Item(4, 20).Value = 0.002452083838914672
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 21.395034
$ws.Cells.Item(5, 8).Value = 64.185102
$ws.Cells.Item(5, 9).Value = 0.03113537258663955
$ws.Cells.Item(5, 10).Value = 0.03113537258663955
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 7.370601
$ws.Cells.Item(5, 14).Value = 22.111803
$ws.Cells.Item(5, 15).Value = 0.2301455394155954
$ws.Cells.Item(5, 16).Value = 0.2301455394155954
$ws.Cells.Item(5, 17).Value = 157.694258995434
$ws.Cells.Item(5, 18).Value = 1419.248330958906
$ws.Cells.Item(5, 19).Value = 0.007165667118857701
$ws.Cells.Item(5, 20).Value = 0.007165667118857699
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 604.8953246666666
$ws.Cells.Item(6, 8).Value = 1814.685974
$ws.Cells.Item(6, 9).Value = 0.8802809712484196
$ws.Cells.Item(6, 10).Value = 0.8802809712484198
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 7.321725333333333
$ws.Cells.Item(6, 14).Value = 21.965176
$ws.Cells.Item(6, 15).Value = 0.2286194065168947
$ws.Cells.Item(6, 16).Value = 0.2286194065168946
$ws.Cells.Item(6, 17).Value = 4428.877422626824
$ws.Cells.Item(6, 18).Value = 39859.89680364142
$ws.Cells.Item(6, 19).Value = 0.2012493132149293
$ws.Cells.Item(6, 20).Value = 0.2012493132149293
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 604.8953246666666
$ws.Cells.Item(7, 8).Value = 1814.685974
$ws.Cells.Item(7, 9).Value = 0.8802809712484196
$ws.Cells.Item(7, 10).Value = 0.8802809712484198
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 14.81128766666667
$ws.Cells.Item(7, 14).Value = 44.433863
$ws.Cells.Item(7, 15).Value = 0.4624794897301531
$ws.Cells.Item(7, 16).Value = 0.462479489730153
$ws.Cells.Item(7, 17).Value = 8959.278661859729
$ws.Cells.Item(7, 18).Value = 80633.50795673757
$ws.Cells.Item(7, 19).Value = 0.4071118944021326
$ws.Cells.Item(7, 20).Value = 0.4071118944021327
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 604.8953246666666
$ws.Cells.Item(8, 8).Value = 1814.685974
$ws.Cells.Item(8, 9).Value = 0.8802809712484196
$ws.Cells.Item(8, 10).Value = 0.8802809712484198
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 2.522212
$ws.Cells.Item(8, 14).Value = 7.566636
$ws.Cells.Item(8, 15).Value = 0.078755564337357
$ws.Cells.Item(8, 16).Value = 0.07875556433735698
$ws.Cells.Item(8, 17).Value = 1525.674246618163
$ws.Cells.Item(8, 18).Value = 13731.06821956346
$ws.Cells.Item(8, 19).Value = 0.06932702466610602
$ws.Cells.Item(8, 20).Value = 0.069327024666106
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 604.8953246666666
$ws.Cells.Item(9, 8).Value = 1814.685974
$ws.Cells.Item(9, 9).Value = 0.8802809712484196
$ws.Cells.Item(9, 10).Value = 0.8802809712484198
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 7.370601
$ws.Cells.Item(9, 14).Value = 22.111803
$ws.Cells.Item(9, 15).Value = 0.2301455394155954
$ws.Cells.Item(9, 16).Value = 0.2301455394155954
$ws.Cells.Item(9, 17).Value = 4458.442084883458
$ws.Cells.Item(9, 18).Value = 40125.97876395112
$ws.Cells.Item(9, 19).Value = 0.2025927389652518
$ws.Cells.Item(9, 20).Value = 0.2025927389652518
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.4378016666666666
$ws.Cells.Item(10, 8).Value = 1.313405
$ws.Cells.Item(10, 9).Value = 0.0006371159779750029
$ws.Cells.Item(10, 10).Value = 0.0006371159779750029
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 7.321725333333333
$ws.Cells.Item(10, 14).Value = 21.965176
$ws.Cells.Item(10, 15).Value = 0.2286194065168947
$ws.Cells.Item(10, 16).Value = 0.2286194065168946
$ws.Cells.Item(10, 17).Value = 3.205463553808888
$ws.Cells.Item(10, 18).Value = 28.84917198428
$ws.Cells.Item(10, 19).Value = 0.0001456570767670761
$ws.Cells.Item(10, 20).Value = 0.0001456570767670761
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.4378016666666666
$ws.Cells.Item(11, 8).Value = 1.313405
$ws.Cells.Item(11, 9).Value = 0.0006371159779750029
$ws.Cells.Item(11, 10).Value = 0.0006371159779750029
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 14.81128766666667
$ws.Cells.Item(11, 14).Value = 44.433863
$ws.Cells.Item(11, 15).Value = 0.4624794897301531
$ws.Cells.Item(11, 16).Value = 0.462479489730153
$ws.Cells.Item(11, 17).Value = 6.484406425946111
$ws.Cells.Item(11, 18).Value = 58.359657833515
$ws.Cells.Item(11, 19).Value = 0.0002946530723928068
$ws.Cells.Item(11, 20).Value = 0.0002946530723928067
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.4378016666666666
$ws.Cells.Item(12, 8).Value = 1.313405
$ws.Cells.Item(12, 9).Value = 0.0006371159779750029
$ws.Cells.Item(12, 10).Value = 0.0006371159779750029
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 2.522212
$ws.Cells.Item(12, 14).Value = 7.566636
$ws.Cells.Item(12, 15).Value = 0.078755564337357
$ws.Cells.Item(12, 16).Value = 0.07875556433735698
$ws.Cells.Item(12, 17).Value = 1.104228617286667
$ws.Cells.Item(12, 18).Value = 9.938057555579999
$ws.Cells.Item(12, 19).Value = 0.00005017642839376847
$ws.Cells.Item(12, 20).Value = 0.00005017642839376845
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.4378016666666666
$ws.Cells.Item(13, 8).Value = 1.313405
$ws.Cells.Item(13, 9).Value = 0.0006371159779750029
$ws.Cells.Item(13, 10).Value = 0.0006371159779750029
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 7.370601
$ws.Cells.Item(13, 14).Value = 22.111803
$ws.Cells.Item(13, 15).Value = 0.2301455394155954
$ws.Cells.Item(13, 16).Value = 0.2301455394155954
$ws.Cells.Item(13, 17).Value = 3.226861402135
$ws.Cells.Item(13, 18).Value = 29.041752619215
$ws.Cells.Item(13, 19).Value = 0.0001466294004213517
$ws.Cells.Item(13, 20).Value = 0.0001466294004213516
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 60.43348966666667
$ws.Cells.Item(14, 8).Value = 181.300469
$ws.Cells.Item(14, 9).Value = 0.08794654018696571
$ws.Cells.Item(14, 10).Value = 0.08794654018696571
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 7.321725333333333
$ws.Cells.Item(14, 14).Value = 21.965176
$ws.Cells.Item(14, 15).Value = 0.2286194065168947
$ws.Cells.Item(14, 16).Value = 0.2286194065168946
$ws.Cells.Item(14, 17).Value = 442.4774122741716
$ws.Cells.Item(14, 18).Value = 3982.296710467544
$ws.Cells.Item(14, 19).Value = 0.02010628582275833
$ws.Cells.Item(14, 20).Value = 0.02010628582275832
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 60.43348966666667
$ws.Cells.Item(15, 8).Value = 181.300469
$ws.Cells.Item(15, 9).Value = 0.08794654018696571
$ws.Cells.Item(15, 10).Value = 0.08794654018696571
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 14.81128766666667
$ws.Cells.Item(15, 14).Value = 44.433863
$ws.Cells.Item(15, 15).Value = 0.4624794897301531
$ws.Cells.Item(15, 16).Value = 0.462479489730153
$ws.Cells.Item(15, 17).Value = 895.0978001535275
$ws.Cells.Item(15, 18).Value = 8055.880201381747
$ws.Cells.Item(15, 19).Value = 0.0406734710292003
$ws.Cells.Item(15, 20).Value = 0.0406734710292003
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 60.43348966666667
$ws.Cells.Item(16, 8).Value = 181.300469
$ws.Cells.Item(16, 9).Value = 0.08794654018696571
$ws.Cells.Item(16, 10).Value = 0.08794654018696571
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 2.522212
$ws.Cells.Item(16, 14).Value = 7.566636
$ws.Cells.Item(16, 15).Value = 0.078755564337357
$ws.Cells.Item(16, 16).Value = 0.07875556433735698
$ws.Cells.Item(16, 17).Value = 152.4260728391427
$ws.Cells.Item(16, 18).Value = 1371.834655552284
$ws.Cells.Item(16, 19).Value = 0.006926279403942531
$ws.Cells.Item(16, 20).Value = 0.006926279403942529
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 60.43348966666667
$ws.Cells.Item(17, 8).Value = 181.300469
$ws.Cells.Item(17, 9).Value = 0.08794654018696571
$ws.Cells.Item(17, 10).Value = 0.08794654018696571
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 7.370601
$ws.Cells.Item(17, 14).Value = 22.111803
$ws.Cells.Item(17, 15).Value = 0.2301455394155954
$ws.Cells.Item(17, 16).Value = 0.2301455394155954
$ws.Cells.Item(17, 17).Value = 445.431139370623
$ws.Cells.Item(17, 18).Value = 4008.880254335606
$ws.Cells.Item(17, 19).Value = 0.02024050393106457
$ws.Cells.Item(17, 20).Value = 0.02024050393106456
